$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "66.879.00"
$ws.Range("E2").Value = "  +2.41%  "
$ws.Range("D3").Value = "3.096.22"
$ws.Range("E3").Value = "  +4.94%  "
$ws.Range("E4").Value = "  +0.09%  "
Set-TextValue "D5" "578.79"
$ws.Range("E5").Value = "  +1.48%  "
Set-TextValue "D6" "172.73"
$ws.Range("E6").Value = "  +7.40%  "
Set-TextValue "D7" "1.00"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "3.092.44"
$ws.Range("E8").Value = "  +4.96%  "
Set-TextValue "D9" "0.523"
$ws.Range("E9").Value = "  +1.42%  "
$ws.Range("E10").Value = "  -2.52%  "
Set-TextValue "D11" "0.155"
$ws.Range("E11").Value = "  +3.79%  "
$ws.Range("E12").Value = "  +5.06%  "
Set-TextValue "D13" "0.0000249"
$ws.Range("E13").Value = "  +2.21%  "
Set-TextValue "D14" "37.13"
$ws.Range("E14").Value = "  +7.65%  "
$ws.Range("E15").Value = "  -0.27%  "
$ws.Range("D16").Value = "3.609.94"
$ws.Range("E16").Value = "  +5.78%  "
$ws.Range("D17").Value = "66.889.34"
$ws.Range("E17").Value = "  +2.47%  "
Set-TextValue "D18" "7.18"
$ws.Range("E18").Value = "  +2.29%  "
$ws.Range("D19").Value = "3.102.41"
$ws.Range("E19").Value = "  +5.19%  "
Set-TextValue "D20" "16.25"
$ws.Range("E20").Value = "  +2.65%  "
Set-TextValue "D21" "483.29"
$ws.Range("E21").Value = "  +8.80%  "
Set-TextValue "D22" "0.714"
$ws.Range("E22").Value = "  +2.48%  "
Set-TextValue "D23" "7.50"
$ws.Range("E23").Value = "  +2.93%  "
Set-TextValue "D24" "84.05"
$ws.Range("E24").Value = "  +2.05%  "
Set-TextValue "D25" "2.34"
$ws.Range("E25").Value = "  +4.74%  "
Set-TextValue "D26" "12.98"
$ws.Range("E26").Value = "  +6.59%  "
$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue "D27" "0.999"
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D28" "9.98"
$ws.Range("E28").Value = "  -0.11%  "
Set-TextValue "D29" "7.98"
$ws.Range("E29").Value = "  +0.51%  "
Set-TextValue "D30" "2.38"
$ws.Range("E30").Value = "  -3.56%  "
Set-TextValue "D31" "2.67"
$ws.Range("E31").Value = "  +3.49%  "
Set-TextValue "D32" "0.0000101"
$ws.Range("E32").Value = "  -0.63%  "
Set-TextValue "D33" "28.71"
$ws.Range("E33").Value = "  +5.86%  "
$ws.Range("E34").Value = "  +2.33%  "
$ws.Range("E35").Value = "  +0.17%  "
Set-TextValue "D36" "1.01"
$ws.Range("E36").Value = "  +3.42%  "
Set-TextValue "D37" "5.87"
$ws.Range("E37").Value = "  +2.56%  "
Set-TextValue "D38" "47.62"
$ws.Range("E38").Value = "  +6.35%  "
Set-TextValue "D39" "2.11"
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D40" "50.16"
$ws.Range("E40").Value = "  +2.03%  "
$ws.Range("B41").Value = "TheGraph"
$ws.Range("C41").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue "D41" "0.315"
$ws.Range("E41").Value = "  +4.84%  "
$ws.Range("E42").Value = "  +0.71%  "
Set-TextValue "D43" "8.66"
$ws.Range("E43").Value = "  +1.21%  "
Set-TextValue "D44" "2.79"
$ws.Range("E44").Value = "  -1.21%  "
Set-TextValue "D45" "0.0359"
$ws.Range("E45").Value = "  +2.50%  "
$ws.Range("D46").Value = "2.816.88"
$ws.Range("E46").Value = "  +5.08%  "
Set-TextValue "D47" "378.55"
$ws.Range("E47").Value = "  -0.70%  "
Set-TextValue "D48" "134.86"
$ws.Range("E48").Value = "  +1.10%  "
$ws.Range("E49").Value = "  -0.02%  "
Set-TextValue "D50" "24.81"
$ws.Range("E50").Value = "  +5.11%  "
Set-TextValue "D51" "2.22"
$ws.Range("E51").Value = "  +2.14%  "
